$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table needs its new entries created in a specific order
# (matching how the original author actually typed them), so the operations
# below are deliberately NOT in simple row order.

# 1) rows 10-13, 16: plain text, keep the default answer style (s=6)
$ws.Range("C10").Value = "The United States"
$ws.Range("C11").Value = "the Andes"
$ws.Range("C12").Value = "Bangkok"
$ws.Range("C13").Value = "The Alps"
$ws.Range("C16").Value = "the Bahamas"

# 2) row 15: rich text ("the" struck through + " Jamaica"), highlighted +
#    quote-prefixed style (s=9). Copy the format from C9 (already s=8 -
#    highlighted) before writing the value so the engine derives / reuses
#    the combined "highlighted + quotePrefix" style.
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "'the Jamaica"
$ws.Range("C15").Characters(1, 3).Font.Strikethrough = $true
$ws.Range("C15").Characters(4, 8).Font.Strikethrough = $false

# 3) row 14: highlighted style (s=8), no quote-prefix needed this time.
$ws.Range("C9").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Kenya(the red sea)"

# 4) rows 19, 22, 24, 27, 29, 34, 36: plain text, default style (s=6)
$ws.Range("C19").Value = "ok"
$ws.Range("C22").Value = "the Regal Cinema"
$ws.Range("C24").Value = "the Museum of Art"
$ws.Range("C27").Value = "the Mississippi or the Nile"
$ws.Range("C29").Value = "the Park Hotel"
$ws.Range("C34").Value = "The west of Ireland"
$ws.Range("C36").Value = "the Panama Canal , the atlantic Ocean , the Pacific Ocean"

# 5) row 31: rich text ("The Rocky Mountains , " + "The " struck through +
#    "North America"), highlighted style (s=8).
$ws.Range("C9").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = "The Rocky Mountains , The North America"
$ws.Range("C31").Characters(23, 4).Font.Strikethrough = $true
$ws.Range("C31").Characters(27, 13).Font.Strikethrough = $false

# 6) rows 28, 26, 21, 20: highlighted style (s=8)
$ws.Range("C9").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = "the london(the national Gallery)"

$ws.Range("C9").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "ok(the netherlands)"

$ws.Range("C9").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = "ok(the south of france)"

$ws.Range("C9").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "ok(the philippines)"

# 7) remaining cells that just reuse already-created shared strings
$ws.Range("C33").Value = "The United States"
$ws.Range("C23").Value = "ok"
$ws.Range("C25").Value = "ok"
$ws.Range("C30").Value = "ok"
$ws.Range("C32").Value = "ok"
$ws.Range("C35").Value = "ok"

# Match the saved selection state
$ws.Range("C19").Select() | Out-Null
